$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.60223533333333
$ws.Range("H2").Value = 70.80670599999999
$ws.Range("I2").Value = 0.8824726436021215
$ws.Range("J2").Value = 0.8824726436021214
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.86952466666667
$ws.Range("N2").Value = 56.608574
$ws.Range("O2").Value = 0.600120086407596
$ws.Range("P2").Value = 0.6001200864075961
$ws.Range("Q2").Value = 445.3629618108048
$ws.Range("R2").Value = 4008.266656297244
$ws.Range("S2").Value = 0.5295895591308448
$ws.Range("T2").Value = 0.5295895591308448

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.60223533333333
$ws.Range("H3").Value = 70.80670599999999
$ws.Range("I3").Value = 0.8824726436021215
$ws.Range("J3").Value = 0.8824726436021214
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.156330666666667
$ws.Range("N3").Value = 6.468992
$ws.Range("O3").Value = 0.0685792233171259
$ws.Range("P3").Value = 0.0685792233171259
$ws.Range("Q3").Value = 50.89422385115022
$ws.Range("R3").Value = 458.048014660352
$ws.Range("S3").Value = 0.06051928849684435
$ws.Range("T3").Value = 0.06051928849684434

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.60223533333333
$ws.Range("H4").Value = 70.80670599999999
$ws.Range("I4").Value = 0.8824726436021215
$ws.Range("J4").Value = 0.8824726436021214
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7389603333333333
$ws.Range("N4").Value = 2.216881
$ws.Range("O4").Value = 0.02350164865971288
$ws.Range("P4").Value = 0.02350164865971288
$ws.Range("Q4").Value = 17.44111568933177
$ws.Range("R4").Value = 156.970041203986
$ws.Range("S4").Value = 0.02073956202174508
$ws.Range("T4").Value = 0.02073956202174508

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.60223533333333
$ws.Range("H5").Value = 70.80670599999999
$ws.Range("I5").Value = 0.8824726436021215
$ws.Range("J5").Value = 0.8824726436021214
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6253503333333333
$ws.Range("N5").Value = 1.876051
$ws.Range("O5").Value = 0.0198884340069237
$ws.Range("P5").Value = 0.0198884340069237
$ws.Range("Q5").Value = 14.75966573311178
$ws.Range("R5").Value = 132.836991598006
$ws.Range("S5").Value = 0.01755099893519629
$ws.Range("T5").Value = 0.01755099893519629

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.60223533333333
$ws.Range("H6").Value = 70.80670599999999
$ws.Range("I6").Value = 0.8824726436021215
$ws.Range("J6").Value = 0.8824726436021214
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.052748666666666
$ws.Range("N6").Value = 27.158246
$ws.Range("O6").Value = 0.2879106076086415
$ws.Range("P6").Value = 0.2879106076086415
$ws.Range("Q6").Value = 213.6651044441862
$ws.Range("R6").Value = 1922.985939997676
$ws.Range("S6").Value = 0.254073235017491
$ws.Range("T6").Value = 0.2540732350174909

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.633202
$ws.Range("H7").Value = 7.899606
$ws.Range("I7").Value = 0.09845375648791208
$ws.Range("J7").Value = 0.09845375648791205
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.86952466666667
$ws.Range("N7").Value = 56.608574
$ws.Range("O7").Value = 0.600120086407596
$ws.Range("P7").Value = 0.6001200864075961
$ws.Range("Q7").Value = 49.68727009131601
$ws.Range("R7").Value = 447.1854308218441
$ws.Range("S7").Value = 0.05908407685067821
$ws.Range("T7").Value = 0.0590840768506782

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.633202
$ws.Range("H8").Value = 7.899606
$ws.Range("I8").Value = 0.09845375648791208
$ws.Range("J8").Value = 0.09845375648791205
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.156330666666667
$ws.Range("N8").Value = 6.468992
$ws.Range("O8").Value = 0.0685792233171259
$ws.Range("P8").Value = 0.0685792233171259
$ws.Range("Q8").Value = 5.678054224128001
$ws.Range("R8").Value = 51.10248801715201
$ws.Range("S8").Value = 0.006751882152594455
$ws.Range("T8").Value = 0.006751882152594453

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.633202
$ws.Range("H9").Value = 7.899606
$ws.Range("I9").Value = 0.09845375648791208
$ws.Range("J9").Value = 0.09845375648791205
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7389603333333333
$ws.Range("N9").Value = 2.216881
$ws.Range("O9").Value = 0.02350164865971288
$ws.Range("P9").Value = 0.02350164865971288
$ws.Range("Q9").Value = 1.945831827654
$ws.Range("R9").Value = 17.512486448886
$ws.Range("S9").Value = 0.002313825594207838
$ws.Range("T9").Value = 0.002313825594207837

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.633202
$ws.Range("H10").Value = 7.899606
$ws.Range("I10").Value = 0.09845375648791208
$ws.Range("J10").Value = 0.09845375648791205
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.6253503333333333
$ws.Range("N10").Value = 1.876051
$ws.Range("O10").Value = 0.0198884340069237
$ws.Range("P10").Value = 0.0198884340069237
$ws.Range("Q10").Value = 1.646673748434
$ws.Range("R10").Value = 14.820063735906
$ws.Range("S10").Value = 0.001958091038643576
$ws.Range("T10").Value = 0.001958091038643575

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.633202
$ws.Range("H11").Value = 7.899606
$ws.Range("I11").Value = 0.09845375648791208
$ws.Range("J11").Value = 0.09845375648791205
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.052748666666666
$ws.Range("N11").Value = 27.158246
$ws.Range("O11").Value = 0.2879106076086415
$ws.Range("P11").Value = 0.2879106076086415
$ws.Range("Q11").Value = 23.837715894564
$ws.Range("R11").Value = 214.539443051076
$ws.Range("S11").Value = 0.028345880851788
$ws.Range("T11").Value = 0.02834588085178799

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.5101343333333334
$ws.Range("H12").Value = 1.530403
$ws.Range("I12").Value = 0.01907359990996641
$ws.Range("J12").Value = 0.0190735999099664
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.86952466666667
$ws.Range("N12").Value = 56.608574
$ws.Range("O12").Value = 0.600120086407596
$ws.Range("P12").Value = 0.6001200864075961
$ws.Range("Q12").Value = 9.625992386146889
$ws.Range("R12").Value = 86.633931475322
$ws.Range("S12").Value = 0.01144645042607296
$ws.Range("T12").Value = 0.01144645042607295

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.5101343333333334
$ws.Range("H13").Value = 1.530403
$ws.Range("I13").Value = 0.01907359990996641
$ws.Range("J13").Value = 0.0190735999099664
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.156330666666667
$ws.Range("N13").Value = 6.468992
$ws.Range("O13").Value = 0.0685792233171259
$ws.Range("P13").Value = 0.0685792233171259
$ws.Range("Q13").Value = 1.100018307086222
$ws.Range("R13").Value = 9.900164763775999
$ws.Range("S13").Value = 0.001308052667687099
$ws.Range("T13").Value = 0.001308052667687098

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5101343333333334
$ws.Range("H14").Value = 1.530403
$ws.Range("I14").Value = 0.01907359990996641
$ws.Range("J14").Value = 0.0190735999099664
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.7389603333333333
$ws.Range("N14").Value = 2.216881
$ws.Range("O14").Value = 0.02350164865971288
$ws.Range("P14").Value = 0.02350164865971288
$ws.Range("Q14").Value = 0.3769690370047778
$ws.Range("R14").Value = 3.392721333043
$ws.Range("S14").Value = 0.0004482610437599618
$ws.Range("T14").Value = 0.0004482610437599617

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5101343333333334
$ws.Range("H15").Value = 1.530403
$ws.Range("I15").Value = 0.01907359990996641
$ws.Range("J15").Value = 0.0190735999099664
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.6253503333333333
$ws.Range("N15").Value = 1.876051
$ws.Range("O15").Value = 0.0198884340069237
$ws.Range("P15").Value = 0.0198884340069237
$ws.Range("Q15").Value = 0.3190126753947778
$ws.Range("R15").Value = 2.871114078553
$ws.Range("S15").Value = 0.0003793440330838328
$ws.Range("T15").Value = 0.0003793440330838326

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5101343333333334
$ws.Range("H16").Value = 1.530403
$ws.Range("I16").Value = 0.01907359990996641
$ws.Range("J16").Value = 0.0190735999099664
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 9.052748666666666
$ws.Range("N16").Value = 27.158246
$ws.Range("O16").Value = 0.2879106076086415
$ws.Range("P16").Value = 0.2879106076086415
$ws.Range("Q16").Value = 4.618117905904223
$ws.Range("R16").Value = 41.56306115313799
$ws.Range("S16").Value = 0.005491491739362558
$ws.Range("T16").Value = 0.005491491739362557
